$d = $word.ActiveDocument

# --- 1. Merge runs that were split by spell-check proofErr markers ---
# These Find/Replace calls target the full visible text of each paragraph;
# since the runtime already concatenates the run text for matching, this
# collapses each split set of runs (separated by w:proofErr) down into the
# single combined run that the diff specifies.

$d.Content.Find.Execute("Added ability to make friend,battle,breed requests", $true, $false, $false, $false, $false, $true, 1, $false, "Added ability to make friend,battle,breed requests", 2) | Out-Null

$d.Content.Find.Execute("Notification reworkings and implementation", $true, $false, $false, $false, $false, $true, 1, $false, "Notification reworkings and implementation", 2) | Out-Null

$d.Content.Find.Execute("Aprox 20hours", $true, $false, $false, $false, $false, $true, 1, $false, "Aprox 20hours", 2) | Out-Null

$d.Content.Find.Execute("See github", $true, $false, $false, $false, $false, $true, 1, $false, "See github", 2) | Out-Null

$d.Content.Find.Execute("Severlet get frends,monster,freindsmonsters", $true, $false, $false, $false, $false, $true, 1, $false, "Severlet get frends,monster,freindsmonsters", 2) | Out-Null

$d.Content.Find.Execute("reworking requests sheat", $true, $false, $false, $false, $false, $true, 1, $false, "reworking requests sheat", 2) | Out-Null

# --- 2. Fill in the empty row (row 10) with the new timesheet entry ---
$t = $d.Tables.Item(1)
$t.Cell(10, 1).Range.Text = "11/01/2013"
$t.Cell(10, 2).Range.Text = "1hour 5am-6am"
$t.Cell(10, 3).Range.Text = "added functianality and continued clean up"
